$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 80, shifting existing rows 80-89 down to 81-90
$ws.Rows.Item(80).Insert()

# Populate the new row 80 with the new record
$ws.Cells.Item(80, 1).Value = 10
$ws.Cells.Item(80, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(80, 3).Value = "La Araucanía"
$ws.Cells.Item(80, 4).Value = 45127
$ws.Cells.Item(80, 5).Value = 9
$ws.Cells.Item(80, 6).Value = 100112026
$ws.Cells.Item(80, 7).Value = "Haba"
$ws.Cells.Item(80, 8).Value = "Sin especificar"
$ws.Cells.Item(80, 9).Value = "Primera"
$ws.Cells.Item(80, 10).Value = 50
$ws.Cells.Item(80, 11).Value = 18000
$ws.Cells.Item(80, 12).Value = 18000
$ws.Cells.Item(80, 13).Value = 18000
$ws.Cells.Item(80, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(80, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(80, 16).Value = 720
$ws.Cells.Item(80, 17).Value = 25
$ws.Cells.Item(80, 18).Value = "Hortaliza"
